# Sprint 6 text updates
# Slide 4 ("Title 1" shape) - bullet list of sprint accomplishments
# Slide 6 ("Text Placeholder 2" shape) - duplicate bullet list on summary slide

function Replace-SubText {
    param($TextRange, $Old, $New)
    $full = $TextRange.Text
    $idx = $full.IndexOf($Old)
    if ($idx -ge 0) {
        $chars = $TextRange.Characters($idx + 1, $Old.Length)
        $chars.Text = $New
    }
}

$p = $ppt.ActivePresentation

# --- Slide 4: Sprint Accomplishments bullet list ---
$slide4 = $p.Slides.Item(4)
$shape4 = $slide4.Shapes.Item(1)
$tr4 = $shape4.TextFrame2.TextRange

Replace-SubText $tr4 "- Added Sign in directory" "- Added Sign in"
Replace-SubText $tr4 "- Emails getting Sent" "- Emails Now getting Sent"
Replace-SubText $tr4 "- Route schedule On Home page" "- Route Modal On Home page"
Replace-SubText $tr4 "- update user experience" "- Fixed Misc. styling"

# --- Slide 6: SUMMARY bullet list (duplicate of the above) ---
$slide6 = $p.Slides.Item(6)
$shape6 = $slide6.Shapes.Item(2)
$tr6 = $shape6.TextFrame2.TextRange

Replace-SubText $tr6 "- Route schedule On Home page" "- Route Modal On Home page"
